$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 14 ("Writing Efficient Python Code"),
# shifting it (and everything below) down by one row.
$ws.Rows.Item(14).Insert()

# Fill in the new row 14 with the newly-added course "Joining Data in SQL"
$ws.Range("A14").Value() = "Joining Data in SQL"
$ws.Range("G14").Value() = 2

# Match the new active selection left behind by the edit
$ws.Range("G15").Select()
